{"js": "// Replace the date line and each three-digit \u00f7 one-digit division problem\n// with the new values, matching the commit's regenerated worksheet output.\nconst replacements = [\n  [\"2024-03-19 Tuesday\", \"2024-03-20 Wednesday\"],\n  [\"913\u00f73=304, 1\", \"365\u00f78=45, 5\"],\n  [\"376\u00f72=188, 0\", \"496\u00f73=165, 1\"],\n  [\"557\u00f76=92, 5\", \"820\u00f74=205, 0\"],\n  [\"951\u00f75=190, 1\", \"184\u00f74=46, 0\"],\n  [\"941\u00f79=104, 5\", \"232\u00f76=38, 4\"],\n  [\"681\u00f78=85, 1\", \"168\u00f77=24, 0\"],\n  [\"231\u00f76=38, 3\", \"780\u00f79=86, 6\"],\n  [\"159\u00f78=19, 7\", \"653\u00f79=72, 5\"],\n  [\"211\u00f77=30, 1\", \"761\u00f76=126, 5\"],\n  [\"591\u00f77=84, 3\", \"397\u00f73=132, 1\"],\n  [\"829\u00f72=414, 1\", \"761\u00f74=190, 1\"],\n  [\"345\u00f76=57, 3\", \"321\u00f74=80, 1\"],\n  [\"930\u00f72=465, 0\", \"203\u00f76=33, 5\"],\n  [\"835\u00f77=119, 2\", \"104\u00f78=13, 0\"],\n  [\"849\u00f75=169, 4\", \"115\u00f73=38, 1\"],\n  [\"956\u00f72=478, 0\", \"889\u00f78=111, 1\"],\n  [\"170\u00f73=56, 2\", \"949\u00f73=316, 1\"],\n  [\"124\u00f73=41, 1\", \"124\u00f76=20, 4\"],\n  [\"209\u00f78=26, 1\", \"909\u00f73=303, 0\"],\n  [\"716\u00f72=358, 0\", \"234\u00f76=39, 0\"],\n  [\"814\u00f77=116, 2\", \"198\u00f76=33, 0\"],\n  [\"959\u00f72=479, 1\", \"456\u00f77=65, 1\"],\n  [\"365\u00f77=52, 1\", \"352\u00f79=39, 1\"],\n  [\"805\u00f79=89, 4\", \"417\u00f76=69, 3\"],\n  [\"510\u00f73=170, 0\", \"749\u00f74=187, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at 9a8706d\n# Replace the date line and each three-digit / one-digit division problem\n# with the regenerated worksheet values.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2024-03-19 Tuesday\"; New = \"2024-03-20 Wednesday\"},\n    @{Old = \"913\u00f73=304, 1\"; New = \"365\u00f78=45, 5\"},\n    @{Old = \"376\u00f72=188, 0\"; New = \"496\u00f73=165, 1\"},\n    @{Old = \"557\u00f76=92, 5\"; New = \"820\u00f74=205, 0\"},\n    @{Old = \"951\u00f75=190, 1\"; New = \"184\u00f74=46, 0\"},\n    @{Old = \"941\u00f79=104, 5\"; New = \"232\u00f76=38, 4\"},\n    @{Old = \"681\u00f78=85, 1\"; New = \"168\u00f77=24, 0\"},\n    @{Old = \"231\u00f76=38, 3\"; New = \"780\u00f79=86, 6\"},\n    @{Old = \"159\u00f78=19, 7\"; New = \"653\u00f79=72, 5\"},\n    @{Old = \"211\u00f77=30, 1\"; New = \"761\u00f76=126, 5\"},\n    @{Old = \"591\u00f77=84, 3\"; New = \"397\u00f73=132, 1\"},\n    @{Old = \"829\u00f72=414, 1\"; New = \"761\u00f74=190, 1\"},\n    @{Old = \"345\u00f76=57, 3\"; New = \"321\u00f74=80, 1\"},\n    @{Old = \"930\u00f72=465, 0\"; New = \"203\u00f76=33, 5\"},\n    @{Old = \"835\u00f77=119, 2\"; New = \"104\u00f78=13, 0\"},\n    @{Old = \"849\u00f75=169, 4\"; New = \"115\u00f73=38, 1\"},\n    @{Old = \"956\u00f72=478, 0\"; New = \"889\u00f78=111, 1\"},\n    @{Old = \"170\u00f73=56, 2\"; New = \"949\u00f73=316, 1\"},\n    @{Old = \"124\u00f73=41, 1\"; New = \"124\u00f76=20, 4\"},\n    @{Old = \"209\u00f78=26, 1\"; New = \"909\u00f73=303, 0\"},\n    @{Old = \"716\u00f72=358, 0\"; New = \"234\u00f76=39, 0\"},\n    @{Old = \"814\u00f77=116, 2\"; New = \"198\u00f76=33, 0\"},\n    @{Old = \"959\u00f72=479, 1\"; New = \"456\u00f77=65, 1\"},\n    @{Old = \"365\u00f77=52, 1\"; New = \"352\u00f79=39, 1\"},\n    @{Old = \"805\u00f79=89, 4\"; New = \"417\u00f76=69, 3\"},\n    @{Old = \"510\u00f73=170, 0\"; New = \"749\u00f74=187, 1\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n}\n"}
